# Apply the "Configuracion HttpGet" update to the PasosExcel workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 0. Highlight the existing "ARCHIVO SQL" row with the light blue accent
#    fill used throughout this update.
# ---------------------------------------------------------------------
$ws.Range("A4:C4").Interior.Color = 15983578

# ---------------------------------------------------------------------
# 1. Add the new "Angular/Ventas/App" section (rows 8-10 first).
# ---------------------------------------------------------------------
$ws.Range("A8").Value = "PROYECTO"
$ws.Range("B8").Value = "Angular/Ventas/App"
$ws.Range("D8").Value = "ng add @angular/material"
$ws.Range("E8").Value = "instalar plantilla de angular"
$ws.Range("A8").Interior.Color = 15983578
$ws.Range("B8").Interior.Color = 15983578
$ws.Range("D8").Interior.Color = 15983578
$ws.Range("E8").Interior.Color = 15983578

$ws.Range("D9").Value = "ng generate component Home"
$ws.Range("E9").Value = "crear componente"

$ws.Range("D10").Value = "ng generate module app-routing --flat --module=app"
$ws.Range("E10").Value = "agregar enrutamiento al proyecto app"

# ---------------------------------------------------------------------
# 2. Split the D2:H2 merge into D2:E2 and F2:G2, and add the
#    "Configuracion" heading in F2:G2.
# ---------------------------------------------------------------------
$ws.Range("D2:H2").UnMerge()
$ws.Range("D2:E2").Merge()
$ws.Range("F2:G2").Merge()

# H2 keeps the same green fill as D2:G2 but is no longer merged with them.
$ws.Range("D2:H2").Interior.Color = 12436269

$ws.Range("F2").Value = "Configuración"
$ws.Range("F2:G2").HorizontalAlignment = -4108

# ---------------------------------------------------------------------
# 3. Add the new "Configuracion" columns (F/G/H) to row 5.
# ---------------------------------------------------------------------
$ws.Range("F5").Value = "Startup.cs"

$corsText = "CORS (Cross Origin Resource Sharing, o bien en español Intercambio de Recursos de Origen Cruzado)"
$ws.Range("H5").Value = $corsText
$ws.Range("H5").Font.Name = "Arial"
$ws.Range("H5").Font.Size = 10
$ws.Range("H5").Font.Bold = $true
$ws.Range("H5").Font.Color = 2367776
$descLen = $corsText.Length - 4
$ws.Range("H5").Characters(5, $descLen).Font.Bold = $false

$ws.Range("G5").Value = "services.AddCors(opciones =>"

# ---------------------------------------------------------------------
# 4. Finish the Angular section with the new service generation row.
# ---------------------------------------------------------------------
$ws.Range("D11").Value = "ng generate service services/apiCliente"
$ws.Range("E11").Value = "crear un servicio; si hay error colocar al inicio npm run"

# ---------------------------------------------------------------------
# 5. Misc sheet-level tweaks to mirror the saved workbook state.
# ---------------------------------------------------------------------
$ws.Range("E13").Select()
$ws.PageSetup.Orientation = 1

$wb.Save()
